# Auto-generated edit script: refresh market-price-driven columns (H,I,J,K,L,M,N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("J3").Value = 58661.75
$ws.Range("L3").Value = 58661.75
$ws.Range("H3").Value = 58661.75
$ws.Range("N3").Value = -58889.75
$ws.Range("M15").Value = -9596
$ws.Range("K15").Value = 9765
$ws.Range("H15").Value = 3255
$ws.Range("I15").Value = 3255
$ws.Range("L102").Value = 58661.75
$ws.Range("J102").Value = 58661.75
$ws.Range("H102").Value = 58661.75
$ws.Range("N102").Value = -65151.75
$ws.Range("K137").Value = 10874.5386
$ws.Range("I137").Value = 3624.8462
$ws.Range("H137").Value = 5922.9375
$ws.Range("M137").Value = -8324.5386

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1629.6154
$ws.Range("K2").Value = 1515.4166
$ws.Range("M2").Value = -1402.4166
$ws.Range("I2").Value = 1515.4166
$ws.Range("K32").Value = 9805368
$ws.Range("H32").Value = 9805368
$ws.Range("I32").Value = 9805368
$ws.Range("M32").Value = -9805081
$ws.Range("H45").Value = 2381.0715
$ws.Range("K45").Value = 2007.3334
$ws.Range("I45").Value = 2007.3334
$ws.Range("M45").Value = -1630.3334
$ws.Range("L61").Value = 119300.11
$ws.Range("M61").Value = -50003164
$ws.Range("I61").Value = 50003376
$ws.Range("H61").Value = 26374076
$ws.Range("N61").Value = -119724.11
$ws.Range("J61").Value = 119300.11
$ws.Range("K61").Value = 50003376
$ws.Range("I74").Value = 16668680
$ws.Range("M74").Value = -16667806
$ws.Range("K74").Value = 16668680
$ws.Range("H74").Value = 7359856
$ws.Range("L74").Value = 10784.263
$ws.Range("J74").Value = 10784.263
$ws.Range("N74").Value = -12532.263
$ws.Range("L77").Value = 53921.315
$ws.Range("I77").Value = 16668680
$ws.Range("N77").Value = -62657.315
$ws.Range("J77").Value = 10784.263
$ws.Range("M77").Value = -83339032
$ws.Range("K77").Value = 83343400
$ws.Range("H77").Value = 7359856
$ws.Range("N82").ClearContents()
$ws.Range("L82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("H82").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("L85").Value = 0
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("H116").Value = 1629.6154
$ws.Range("I116").Value = 1515.4166
$ws.Range("K116").Value = 1515.4166
$ws.Range("M116").Value = 778.5834
$ws.Range("M122").Value = -8059.999899999999
$ws.Range("H122").Value = 4438.75
$ws.Range("I122").Value = 3503.3333
$ws.Range("K122").Value = 10509.9999
$ws.Range("M136").Value = -150007578
$ws.Range("J136").Value = 119300.11
$ws.Range("I136").Value = 50003376
$ws.Range("H136").Value = 26374076
$ws.Range("L136").Value = 357900.33
$ws.Range("K136").Value = 150010128
$ws.Range("N136").Value = -363000.33

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("K3").Value = 1515.4166
$ws.Range("H3").Value = 1629.6154
$ws.Range("M3").Value = -1401.4166
$ws.Range("I3").Value = 1515.4166
$ws.Range("M86").Value = -3.894800000000032
$ws.Range("K86").Value = 1126.8948
$ws.Range("N86").Value = -3699.5
$ws.Range("J86").Value = 1453.5
$ws.Range("L86").Value = 1453.5
$ws.Range("H86").Value = 1158
$ws.Range("I86").Value = 1126.8948
$ws.Range("I89").Value = 1126.8948
$ws.Range("J89").Value = 1453.5
$ws.Range("K89").Value = 5634.474
$ws.Range("M89").Value = -18.47400000000016
$ws.Range("N89").Value = -18499.5
$ws.Range("H89").Value = 1158
$ws.Range("L89").Value = 7267.5
$ws.Range("K107").Value = 2550.4443
$ws.Range("L107").Value = 2013
$ws.Range("I107").Value = 2550.4443
$ws.Range("H107").Value = 2496.7
$ws.Range("J107").Value = 2013
$ws.Range("N107").Value = -5853
$ws.Range("M107").Value = -630.4443000000001
$ws.Range("I134").Value = 3444.2964
$ws.Range("K134").Value = 10332.8892
$ws.Range("M134").Value = -7797.889200000001
$ws.Range("H134").Value = 36056.676

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("K16").Value = 974.1818
$ws.Range("H16").Value = 1294
$ws.Range("I16").Value = 974.1818
$ws.Range("M16").Value = -687.1818
$ws.Range("K31").Value = 1665.8462
$ws.Range("M31").Value = -1370.8462
$ws.Range("I31").Value = 1665.8462
$ws.Range("H31").Value = 780035.0600000001
$ws.Range("M34").Value = -1463.8462
$ws.Range("I34").Value = 1665.8462
$ws.Range("H34").Value = 780035.0600000001
$ws.Range("K34").Value = 1665.8462
$ws.Range("J108").Value = 65056.875
$ws.Range("L108").Value = 65056.875
$ws.Range("H108").Value = 65056.875
$ws.Range("N108").Value = -72736.875
$ws.Range("H113").Value = 1294
$ws.Range("M113").Value = 1195.8182
$ws.Range("I113").Value = 974.1818
$ws.Range("K113").Value = 974.1818

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I4").Value = 131937544
$ws.Range("H4").Value = 113815830
$ws.Range("M4").Value = -395812520
$ws.Range("K4").Value = 395812632
$ws.Range("M8").ClearContents()
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("N16").Value = -22097.5
$ws.Range("K16").Value = 814.5
$ws.Range("H16").Value = 2265.5
$ws.Range("I16").Value = 271.5
$ws.Range("J16").Value = 7250.5
$ws.Range("L16").Value = 21751.5
$ws.Range("M16").Value = -641.5
$ws.Range("I23").Value = 263.8
$ws.Range("K23").Value = 791.4000000000001
$ws.Range("J23").Value = 293.8889
$ws.Range("H23").Value = 283.14285
$ws.Range("L23").Value = 881.6667
$ws.Range("M23").Value = -556.4000000000001
$ws.Range("N23").Value = -1351.6667
$ws.Range("N35").ClearContents()
$ws.Range("M35").Value = -3912
$ws.Range("I35").Value = 1400
$ws.Range("L35").Value = 0
$ws.Range("H35").Value = 1400
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 4200
$ws.Range("K68").Value = 5925
$ws.Range("M68").Value = -5114
$ws.Range("H68").Value = 2560.9
$ws.Range("I68").Value = 1975
$ws.Range("I71").Value = 1975
$ws.Range("M71").Value = -13719
$ws.Range("K71").Value = 17775
$ws.Range("H71").Value = 2560.9
$ws.Range("M80").Value = -5811
$ws.Range("K80").Value = 6747
$ws.Range("J80").Value = 3750
$ws.Range("N80").Value = -13122
$ws.Range("I80").Value = 2249
$ws.Range("L80").Value = 11250
$ws.Range("H80").Value = 3562.375
$ws.Range("M83").Value = -15561
$ws.Range("H83").Value = 3562.375
$ws.Range("K83").Value = 20241
$ws.Range("J83").Value = 3750
$ws.Range("I83").Value = 2249
$ws.Range("L83").Value = 33750
$ws.Range("N83").Value = -43110
$ws.Range("J109").Value = 2300
$ws.Range("N109").Value = -8980
$ws.Range("L109").Value = 6900
$ws.Range("H109").Value = 2543.3333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4287.75
$ws.Range("M113").Value = -815
$ws.Range("I113").Value = 2985
$ws.Range("K113").Value = 2985
$ws.Range("M132").Value = -214290902
$ws.Range("H132").Value = 71431144
$ws.Range("I132").Value = 71431144
$ws.Range("K132").Value = 214293432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("K22").Value = 1034.6
$ws.Range("M22").Value = -739.5999999999999
$ws.Range("I22").Value = 1034.6
$ws.Range("H22").Value = 1028.1177
$ws.Range("M27").Value = -927.5999999999999
$ws.Range("K27").Value = 1034.6
$ws.Range("I27").Value = 1034.6
$ws.Range("H27").Value = 1028.1177
$ws.Range("K40").Value = 4360.4
$ws.Range("I40").Value = 4360.4
$ws.Range("M40").Value = -4224.4
$ws.Range("H40").Value = 5254.636
$ws.Range("M136").Value = -433438.26
$ws.Range("I136").Value = 145329.42
$ws.Range("H136").Value = 146487.4
$ws.Range("K136").Value = 435988.26

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("J75").Value = 12543301
$ws.Range("N75").Value = -12545173
$ws.Range("L75").Value = 12543301
$ws.Range("H75").Value = 12543301
$ws.Range("N78").Value = -37639263
$ws.Range("L78").Value = 37629903
$ws.Range("J78").Value = 12543301
$ws.Range("H78").Value = 12543301
$ws.Range("M136").Value = -1623.800099999999
$ws.Range("I136").Value = 1391.2667
$ws.Range("H136").Value = 1462.7646
$ws.Range("K136").Value = 4173.800099999999
